# Auto-generated edit script: updates the Goshen_A transition-matrix values
# per the commit "added more games, sped up simulate game logic, and drafted optimization logic"
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 0.1
$ws.Range("C2").Value = 0.65
$ws.Range("P2").Value = 0.1
$ws.Range("S2").Value = 0.15
$ws.Range("J3").Value = 0.07142857142857142
$ws.Range("P3").Value = 0.7142857142857143
$ws.Range("S3").Value = 0.2142857142857143
$ws.Range("P4").Value = 0.6666666666666666
$ws.Range("S4").Value = 0.3333333333333333
$ws.Range("D6").Value = 0.05
$ws.Range("F6").Value = 0.15
$ws.Range("J6").Value = 0.35
$ws.Range("Q6").Value = 0.1
$ws.Range("S6").Value = 0.35
$ws.Range("B7").Value = 0.04761904761904762
$ws.Range("J7").Value = 0.1428571428571428
$ws.Range("O7").Value = 0.09523809523809523
$ws.Range("Q7").Value = 0.09523809523809523
$ws.Range("R7").Value = 0.04761904761904762
$ws.Range("S7").Value = 0.5714285714285714
$ws.Range("B8").Value = 0.1333333333333333
$ws.Range("D8").Value = 0.04444444444444445
$ws.Range("F8").Value = 0.06666666666666667
$ws.Range("J8").Value = 0.02222222222222222
$ws.Range("O8").Value = 0.02222222222222222
$ws.Range("Q8").Value = 0.08888888888888889
$ws.Range("R8").Value = 0.08888888888888889
$ws.Range("S8").Value = 0.5333333333333333
$ws.Range("B9").Value = 0.1538461538461539
$ws.Range("F9").Value = 0.03846153846153846
$ws.Range("J9").Value = 0.07692307692307693
$ws.Range("Q9").Value = 0.1923076923076923
$ws.Range("R9").Value = 0.07692307692307693
$ws.Range("S9").Value = 0.4615384615384616
$ws.Range("B10").Value = 0.109375
$ws.Range("F10").Value = 0.046875
$ws.Range("J10").Value = 0.125
$ws.Range("O10").Value = 0.03125
$ws.Range("Q10").Value = 0.109375
$ws.Range("R10").Value = 0.03125
$ws.Range("S10").Value = 0.546875
$ws.Range("G11").Value = 0.3333333333333333
$ws.Range("J11").Value = 0.03703703703703703
$ws.Range("L11").Value = 0.2962962962962963
$ws.Range("G12").Value = 0.625
$ws.Range("J12").Value = 0.125
$ws.Range("S12").Value = 0.25
$ws.Range("G13").Value = 0.875
$ws.Range("J13").Value = 0.125
$ws.Range("F15").Value = 0.07142857142857142
$ws.Range("I15").Value = 0.07142857142857142
$ws.Range("J15").Value = 0.07142857142857142
$ws.Range("K15").Value = 0.2142857142857143
$ws.Range("M15").Value = 0.07142857142857142
$ws.Range("O15").Value = 0.1428571428571428
$ws.Range("S15").Value = 0.3571428571428572
$ws.Range("H16").Value = 0.2857142857142857
$ws.Range("I16").Value = 0.07142857142857142
$ws.Range("J16").Value = 0.1428571428571428
$ws.Range("K16").Value = 0.07142857142857142
$ws.Range("M16").Value = 0.1428571428571428
$ws.Range("O16").Value = 0.07142857142857142
$ws.Range("S16").Value = 0.2142857142857143
$ws.Range("F17").Value = 0.1052631578947368
$ws.Range("H17").Value = 0.3157894736842105
$ws.Range("I17").Value = 0.1578947368421053
$ws.Range("J17").Value = 0.1578947368421053
$ws.Range("K17").Value = 0.1052631578947368
$ws.Range("O17").Value = 0.05263157894736842
$ws.Range("S17").Value = 0.1052631578947368
$ws.Range("F18").Value = 0.1111111111111111
$ws.Range("H18").Value = 0.1111111111111111
$ws.Range("I18").Value = 0.1111111111111111
$ws.Range("J18").Value = 0.2222222222222222
$ws.Range("K18").Value = 0.1111111111111111
$ws.Range("M18").Value = 0.1111111111111111
$ws.Range("S18").Value = 0.2222222222222222
$ws.Range("F19").Value = 0.03076923076923077
$ws.Range("H19").Value = 0.2615384615384616
$ws.Range("I19").Value = 0.1538461538461539
$ws.Range("J19").Value = 0.2615384615384616
$ws.Range("K19").Value = 0.07692307692307693
$ws.Range("M19").Value = 0.03076923076923077
$ws.Range("O19").Value = 0.03076923076923077
$ws.Range("S19").Value = 0.1538461538461539
